$wb = $excel.ActiveWorkbook

# --- Update the "Additional Body Text" for the "Clean/Dirty Install Smoke
# Tests" row on the "issues" sheet: remove the outdated/obsolete
# clean-install instruction (the "Find `HKEY_CURRENT_USER ..." line) and
# tidy up the blank separator lines. ---
$issues = $wb.Worksheets.Item("issues")

$newText = @'
### Dirty install
* Make sure that you have several versions of Mantid installed
 * Last release
 * A nightly
 * If possible an old release
* Install the latest version of the new Mantid
- [ ] Check that Mantid boots up correctly
### Clean install
* Remove all existing Mantid versions and associated files
* Windows:
 * Uninstall the program
 * Clear shortcuts from desktop
 * Clean out the registry
  * Load regedit (Command Prompt > regedit)
**On macOS** :
 * Remove the application
 * Remove the `~/.mantid directory`
 * Remove `~/Library/Preferences/org.mantidproject.MantidPlot.plist`
**On Linux** :
 * Remove the package: `/opt/Mantid`
 * Remove `~/.config/Mantid`
 * Remove `~/.mantid/`
* Re-install the latest version of the new Mantid
- [ ] Check that Mantid boots up correctly

'@

$issues.Range("C2").Value = $newText

# Assigning the (shorter) text can make Excel auto-fit the row to a new
# height; restore the original explicit row height so only the cell
# content/shared-string reference changes.
$issues.Rows.Item(2).RowHeight = 135

# --- Make "issues" the active sheet/tab, with cell G4 selected (matches
# the saved view state) ---
$issues.Activate()
$issues.Range("G4").Select()
